# Update cryptocurrency price (D) and volume-change (E) columns
# with the latest scraped figures. D-column values are stored as
# plain text (matching the source feed formatting, incl. trailing
# zeros and thousands-dot separators), so we force a Text number
# format while assigning, then restore the default "Normal" style
# so no stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.632.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.225.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.10%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.224.16'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.783.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.739.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.236.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000160'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '414.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.206'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.495'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  -4.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.02'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.771.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.15%  '
$ws.Range("E44").Value = '  -6.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0631'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '302.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.00%  '
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("E51").Value = '  -2.68%  '
